# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows (343-345) into the Limon price list sheet,
# pushing the existing rows 343:417 down to 346:420.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 343; everything currently at
# rows 343:417 (and their formatting) shifts down to rows 346:420.
$ws.Range("A343:T345").Insert()

# Constant columns shared by every record in this sheet.
$mercadoId = 7
$mercado   = "Terminal Hortofrutícola Agro Chillán"
$region    = "Ñuble"
$codreg    = 16
$tipo      = "Fruta"
$productoId = 100102
$producto   = "Cítricos"
$categoriaId = 100102003
$categoria   = "Limón"
$variedad    = "Sin especificar"

# New row 343: 1a amarillo
$ws.Cells.Item(343, 1).Value = $mercadoId
$ws.Cells.Item(343, 2).Value = $mercado
$ws.Cells.Item(343, 3).Value = $region
$ws.Range("D343").Value = 44476
$ws.Cells.Item(343, 5).Value = $codreg
$ws.Cells.Item(343, 6).Value = $tipo
$ws.Cells.Item(343, 7).Value = $productoId
$ws.Cells.Item(343, 8).Value = $producto
$ws.Cells.Item(343, 9).Value = $categoriaId
$ws.Cells.Item(343, 10).Value = $categoria
$ws.Cells.Item(343, 11).Value = $variedad
$ws.Cells.Item(343, 12).Value = "1a amarillo"
$ws.Cells.Item(343, 13).Value = 240
$ws.Cells.Item(343, 14).Value = 4000
$ws.Cells.Item(343, 15).Value = 4500
$ws.Cells.Item(343, 16).Value = 4250
$ws.Cells.Item(343, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(343, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(343, 19).Value = 266
$ws.Cells.Item(343, 20).Value = 16

# New row 344: 2a amarillo
$ws.Cells.Item(344, 1).Value = $mercadoId
$ws.Cells.Item(344, 2).Value = $mercado
$ws.Cells.Item(344, 3).Value = $region
$ws.Range("D344").Value = 44476
$ws.Cells.Item(344, 5).Value = $codreg
$ws.Cells.Item(344, 6).Value = $tipo
$ws.Cells.Item(344, 7).Value = $productoId
$ws.Cells.Item(344, 8).Value = $producto
$ws.Cells.Item(344, 9).Value = $categoriaId
$ws.Cells.Item(344, 10).Value = $categoria
$ws.Cells.Item(344, 11).Value = $variedad
$ws.Cells.Item(344, 12).Value = "2a amarillo"
$ws.Cells.Item(344, 13).Value = 240
$ws.Cells.Item(344, 14).Value = 3000
$ws.Cells.Item(344, 15).Value = 3500
$ws.Cells.Item(344, 16).Value = 3250
$ws.Cells.Item(344, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(344, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(344, 19).Value = 203
$ws.Cells.Item(344, 20).Value = 16

# New row 345: 3a amarillo
$ws.Cells.Item(345, 1).Value = $mercadoId
$ws.Cells.Item(345, 2).Value = $mercado
$ws.Cells.Item(345, 3).Value = $region
$ws.Range("D345").Value = 44476
$ws.Cells.Item(345, 5).Value = $codreg
$ws.Cells.Item(345, 6).Value = $tipo
$ws.Cells.Item(345, 7).Value = $productoId
$ws.Cells.Item(345, 8).Value = $producto
$ws.Cells.Item(345, 9).Value = $categoriaId
$ws.Cells.Item(345, 10).Value = $categoria
$ws.Cells.Item(345, 11).Value = $variedad
$ws.Cells.Item(345, 12).Value = "3a amarillo"
$ws.Cells.Item(345, 13).Value = 60
$ws.Cells.Item(345, 14).Value = 2800
$ws.Cells.Item(345, 15).Value = 2800
$ws.Cells.Item(345, 16).Value = 2800
$ws.Cells.Item(345, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(345, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(345, 19).Value = 175
$ws.Cells.Item(345, 20).Value = 16
